$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared-string rich-text runs) ---
$ws.Range("A8").Characters(21, 2).Text = "22"
$ws.Range("C9").Characters(27, 9).Text = "5/29/2023"
$ws.Range("C9").Characters(47, 9).Text = "6/4/2023"

# --- Weekly crime-statistics table updates (rows 14-30) ---
$ws.Range("C14").Value = 5
$ws.Range("D14").Value = 5
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 28
$ws.Range("G14").Value = 38
$ws.Range("H14").Value = -26.315789473684
$ws.Range("I14").Value = 161
$ws.Range("J14").Value = 185
$ws.Range("K14").Value = -12.972972972973
$ws.Range("L14").Value = -16.580310880829
$ws.Range("M14").Value = -16.580310880829
$ws.Range("N14").Value = -79.773869346733
$ws.Range("C15").Value = 23
$ws.Range("D15").Value = 27
$ws.Range("E15").Value = -14.814814814814
$ws.Range("F15").Value = 120
$ws.Range("G15").Value = 142
$ws.Range("H15").Value = -15.492957746478
$ws.Range("I15").Value = 631
$ws.Range("J15").Value = 680
$ws.Range("K15").Value = -7.205882352941
$ws.Range("L15").Value = 6.949152542372
$ws.Range("M15").Value = 21.113243761996
$ws.Range("N15").Value = -53.840526700804
$ws.Range("C16").Value = 314
$ws.Range("D16").Value = 369
$ws.Range("E16").Value = -14.905149051490
$ws.Range("F16").Value = 1224
$ws.Range("G16").Value = 1373
$ws.Range("H16").Value = -10.852148579752
$ws.Range("I16").Value = 6541
$ws.Range("J16").Value = 6800
$ws.Range("K16").Value = -3.808823529411
$ws.Range("L16").Value = 34.921617161716
$ws.Range("M16").Value = -13.638764193292
$ws.Range("N16").Value = -81.335425881010
$ws.Range("C17").Value = 541
$ws.Range("D17").Value = 566
$ws.Range("E17").Value = -4.416961130742
$ws.Range("F17").Value = 2166
$ws.Range("G17").Value = 2185
$ws.Range("H17").Value = -0.869565217391
$ws.Range("I17").Value = 11101
$ws.Range("J17").Value = 10379
$ws.Range("K17").Value = 6.956354176702
$ws.Range("L17").Value = 29.231664726426
$ws.Range("M17").Value = 60.141373341027
$ws.Range("N17").Value = -32.331606217616
$ws.Range("C18").Value = 247
$ws.Range("D18").Value = 314
$ws.Range("E18").Value = -21.337579617834
$ws.Range("F18").Value = 1005
$ws.Range("G18").Value = 1162
$ws.Range("H18").Value = -13.511187607573
$ws.Range("I18").Value = 5963
$ws.Range("J18").Value = 6502
$ws.Range("K18").Value = -8.289756997846
$ws.Range("L18").Value = 22.720724428894
$ws.Range("M18").Value = -20.034866568325
$ws.Range("N18").Value = -85.661039773000
$ws.Range("C19").Value = 972
$ws.Range("D19").Value = 945
$ws.Range("E19").Value = 2.857142857142
$ws.Range("F19").Value = 3837
$ws.Range("G19").Value = 3785
$ws.Range("H19").Value = 1.373844121532
$ws.Range("I19").Value = 20518
$ws.Range("J19").Value = 20697
$ws.Range("K19").Value = -0.864859641493
$ws.Range("L19").Value = 50.281989306379
$ws.Range("M19").Value = 38.382680245498
$ws.Range("N19").Value = -39.981278886093
$ws.Range("C20").Value = 250
$ws.Range("D20").Value = 232
$ws.Range("E20").Value = 7.758620689655
$ws.Range("F20").Value = 1185
$ws.Range("G20").Value = 949
$ws.Range("H20").Value = 24.868282402529
$ws.Range("I20").Value = 6304
$ws.Range("J20").Value = 5379
$ws.Range("K20").Value = 17.196504926566
$ws.Range("L20").Value = 76.780706674144
$ws.Range("M20").Value = 50.417561441183
$ws.Range("N20").Value = -86.771587451474
$ws.Range("C21").Value = 2352
$ws.Range("D21").Value = 2458
$ws.Range("E21").Value = -4.312449145646
$ws.Range("F21").Value = 9565
$ws.Range("G21").Value = 9634
$ws.Range("H21").Value = -0.716213410836
$ws.Range("I21").Value = 51219
$ws.Range("J21").Value = 50622
$ws.Range("K21").Value = 1.179329145430
$ws.Range("L21").Value = 41.103060690377
$ws.Range("M21").Value = 22.842067394172
$ws.Range("N21").Value = -71.069249887031
$ws.Range("C22").Value = 39
$ws.Range("D22").Value = 48
$ws.Range("E22").Value = -18.75
$ws.Range("F22").Value = 180
$ws.Range("G22").Value = 193
$ws.Range("H22").Value = -6.735751295336
$ws.Range("I22").Value = 916
$ws.Range("J22").Value = 996
$ws.Range("K22").Value = -8.032128514056
$ws.Range("L22").Value = 43.348982785602
$ws.Range("M22").Value = 2.805836139169
$ws.Range("C23").Value = 130
$ws.Range("D23").Value = 135
$ws.Range("E23").Value = -3.703703703703
$ws.Range("F23").Value = 463
$ws.Range("G23").Value = 489
$ws.Range("H23").Value = -5.316973415132
$ws.Range("I23").Value = 2563
$ws.Range("J23").Value = 2415
$ws.Range("K23").Value = 6.128364389233
$ws.Range("L23").Value = 17.461044912923
$ws.Range("M23").Value = 58.014796547472
$ws.Range("C24").Value = 2194
$ws.Range("D24").Value = 2287
$ws.Range("E24").Value = -4.066462614779
$ws.Range("F24").Value = 8664
$ws.Range("G24").Value = 9089
$ws.Range("H24").Value = -4.675981956210
$ws.Range("I24").Value = 45443
$ws.Range("J24").Value = 46093
$ws.Range("K24").Value = -1.410192437029
$ws.Range("L24").Value = 40.429542645241
$ws.Range("M24").Value = 40.703470910610
$ws.Range("C25").Value = 950
$ws.Range("D25").Value = 937
$ws.Range("E25").Value = 1.387406616862
$ws.Range("F25").Value = 3694
$ws.Range("G25").Value = 3462
$ws.Range("H25").Value = 6.701328711727
$ws.Range("I25").Value = 18067
$ws.Range("J25").Value = 17078
$ws.Range("K25").Value = 5.791076238435
$ws.Range("L25").Value = 34.707724425887
$ws.Range("M25").Value = -5.254601709580
$ws.Range("C26").Value = 44
$ws.Range("D26").Value = 47
$ws.Range("E26").Value = -6.382978723404
$ws.Range("F26").Value = 215
$ws.Range("G26").Value = 242
$ws.Range("H26").Value = -11.157024793388
$ws.Range("I26").Value = 1046
$ws.Range("J26").Value = 1123
$ws.Range("K26").Value = -6.856634016028
$ws.Range("L26").Value = 5.870445344129
$ws.Range("C27").Value = 110
$ws.Range("D27").Value = 112
$ws.Range("E27").Value = -1.785714285714
$ws.Range("F27").Value = 484
$ws.Range("G27").Value = 471
$ws.Range("H27").Value = 2.760084925690
$ws.Range("I27").Value = 2205
$ws.Range("J27").Value = 2100
$ws.Range("K27").Value = 5
$ws.Range("L27").Value = 19.771863117870
$ws.Range("C28").Value = 25
$ws.Range("D28").Value = 38
$ws.Range("E28").Value = -34.210526315789
$ws.Range("F28").Value = 92
$ws.Range("G28").Value = 126
$ws.Range("H28").Value = -26.984126984127
$ws.Range("I28").Value = 466
$ws.Range("J28").Value = 619
$ws.Range("K28").Value = -24.717285945072
$ws.Range("L28").Value = -30.134932533733
$ws.Range("M28").Value = -27.975270479134
$ws.Range("N28").Value = -80.304311073541
$ws.Range("C29").Value = 24
$ws.Range("D29").Value = 33
$ws.Range("E29").Value = -27.272727272727
$ws.Range("F29").Value = 84
$ws.Range("G29").Value = 108
$ws.Range("H29").Value = -22.222222222222
$ws.Range("I29").Value = 398
$ws.Range("J29").Value = 527
$ws.Range("K29").Value = -24.478178368121
$ws.Range("L29").Value = -32.081911262798
$ws.Range("M29").Value = -25.047080979284
$ws.Range("N29").Value = -81.427904806346
$ws.Range("C30").Value = 3
$ws.Range("D30").Value = 9
$ws.Range("E30").Value = -66.666666666666
$ws.Range("F30").Value = 34
$ws.Range("G30").Value = 47
$ws.Range("H30").Value = -27.659574468085
$ws.Range("I30").Value = 204
$ws.Range("J30").Value = 291
$ws.Range("K30").Value = -29.896907216494
$ws.Range("L30").Value = -16.734693877551
